$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 was a shared string "3 tai 4"; it becomes a plain number (5) like the
# other hour-count cells below it.
$ws.Range("B4").Value = 5

# Apply an integer number format ("0") to the whole hour-count column for
# the existing rows and a block of blank rows below (this mints cellXfs
# entry #2 = numFmtId 1 in styles.xml, reused by every B cell touched).
$ws.Range("B2:B34").NumberFormat = "0"

# New row 7: 12-14.4 / 10 / Tutkiskelin omaa koodiani ja opettelin pythonia
$ws.Range("A7").Value = "12-14.4"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = "Tutkiskelin omaa koodiani ja opettelin pythonia"

# New row 8: 15.4.2013 (date, same d-mmm style as the other date cells) /
# 3 / tapaaminen ja lisakoodausta
$ws.Range("A8").Value = 41379
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "tapaaminen ja lisäkoodausta"

# Touching PageSetup mints a printerSettings relationship + <pageSetup/>.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection now sits on B8, matching where data entry left off.
$ws.Range("B8").Select()
